# Update the F-column "time_taken" timestamps on the "data" sheet to reflect
# a fresh query run, then add a new "metadata" worksheet summarising the panel
# fetch (panel name/id/version/request info), mirroring the author's commit:
# "Refined metadata to be additional tab".

$wb = $excel.ActiveWorkbook
$dataWs = $wb.Worksheets.Item("data")

# ---------------------------------------------------------------------------
# 1) Refresh the per-row "time_taken" timestamps in column F of "data".
# ---------------------------------------------------------------------------
$timeUpdates = @(
    @(2, "2021-10-05 14:20:26.041628"),
    @(3, "2021-10-05 14:20:26.041651"),
    @(4, "2021-10-05 14:20:26.041657"),
    @(5, "2021-10-05 14:20:26.041660"),
    @(6, "2021-10-05 14:20:26.041663"),
    @(7, "2021-10-05 14:20:26.041665"),
    @(8, "2021-10-05 14:20:26.041668"),
    @(9, "2021-10-05 14:20:26.041671"),
    @(10, "2021-10-05 14:20:26.041673"),
    @(11, "2021-10-05 14:20:26.041676"),
    @(12, "2021-10-05 14:20:26.041679"),
    @(13, "2021-10-05 14:20:26.041681"),
    @(14, "2021-10-05 14:20:26.041684"),
    @(15, "2021-10-05 14:20:26.041686"),
    @(16, "2021-10-05 14:20:26.041689"),
    @(17, "2021-10-05 14:20:26.041691"),
    @(18, "2021-10-05 14:20:26.041694"),
    @(19, "2021-10-05 14:20:26.041697"),
    @(20, "2021-10-05 14:20:26.041700"),
    @(21, "2021-10-05 14:20:26.041702"),
    @(22, "2021-10-05 14:20:26.041705"),
    @(23, "2021-10-05 14:20:26.041708"),
    @(24, "2021-10-05 14:20:26.041710"),
    @(25, "2021-10-05 14:20:26.041713"),
    @(26, "2021-10-05 14:20:26.041716"),
    @(27, "2021-10-05 14:20:26.041718"),
    @(28, "2021-10-05 14:20:26.041721"),
    @(29, "2021-10-05 14:20:26.041724"),
    @(30, "2021-10-05 14:20:26.041726"),
    @(31, "2021-10-05 14:20:26.041729"),
    @(32, "2021-10-05 14:20:26.041732"),
    @(33, "2021-10-05 14:20:26.041734"),
    @(34, "2021-10-05 14:20:26.041737"),
    @(35, "2021-10-05 14:20:26.041740"),
    @(36, "2021-10-05 14:20:26.041742"),
    @(37, "2021-10-05 14:20:26.041745"),
    @(38, "2021-10-05 14:20:26.041747"),
    @(39, "2021-10-05 14:20:26.041750"),
    @(40, "2021-10-05 14:20:26.041752"),
    @(41, "2021-10-05 14:20:26.041755"),
    @(42, "2021-10-05 14:20:26.041758"),
    @(43, "2021-10-05 14:20:26.041761"),
    @(44, "2021-10-05 14:20:26.041763"),
    @(45, "2021-10-05 14:20:26.041766"),
    @(46, "2021-10-05 14:20:26.041768"),
    @(47, "2021-10-05 14:20:26.041771"),
    @(48, "2021-10-05 14:20:26.041773"),
    @(49, "2021-10-05 14:20:26.041776"),
    @(50, "2021-10-05 14:20:26.041778"),
    @(51, "2021-10-05 14:20:26.041781"),
    @(52, "2021-10-05 14:20:26.041784"),
    @(53, "2021-10-05 14:20:26.041786"),
    @(54, "2021-10-05 14:20:26.041789"),
    @(55, "2021-10-05 14:20:26.041792"),
    @(56, "2021-10-05 14:20:26.041794"),
    @(57, "2021-10-05 14:20:26.041797"),
    @(58, "2021-10-05 14:20:26.041799"),
    @(59, "2021-10-05 14:20:26.041802"),
    @(60, "2021-10-05 14:20:26.041804"),
    @(61, "2021-10-05 14:20:26.041807"),
    @(62, "2021-10-05 14:20:26.041809"),
    @(63, "2021-10-05 14:20:26.041812"),
    @(64, "2021-10-05 14:20:26.041814"),
    @(65, "2021-10-05 14:20:26.041817"),
    @(66, "2021-10-05 14:20:26.041821"),
    @(67, "2021-10-05 14:20:26.041824"),
    @(68, "2021-10-05 14:20:26.041827"),
    @(69, "2021-10-05 14:20:26.041830"),
    @(70, "2021-10-05 14:20:26.041832"),
    @(71, "2021-10-05 14:20:26.041835"),
    @(72, "2021-10-05 14:20:26.041838"),
    @(73, "2021-10-05 14:20:26.041840"),
    @(74, "2021-10-05 14:20:26.041843"),
    @(75, "2021-10-05 14:20:26.041845"),
    @(76, "2021-10-05 14:20:26.041848"),
    @(77, "2021-10-05 14:20:26.041851"),
    @(78, "2021-10-05 14:20:26.041855"),
    @(79, "2021-10-05 14:20:26.041858"),
    @(80, "2021-10-05 14:20:26.041861"),
    @(81, "2021-10-05 14:20:26.041864"),
    @(82, "2021-10-05 14:20:26.041866"),
    @(83, "2021-10-05 14:20:26.041869"),
    @(84, "2021-10-05 14:20:26.041871"),
    @(85, "2021-10-05 14:20:26.041874"),
    @(86, "2021-10-05 14:20:26.041876"),
    @(87, "2021-10-05 14:20:26.041879"),
    @(88, "2021-10-05 14:20:26.041882"),
    @(89, "2021-10-05 14:20:26.041884"),
    @(90, "2021-10-05 14:20:26.041887")
)

foreach ($pair in $timeUpdates) {
    $row = $pair[0]
    $newValue = $pair[1]
    $dataWs.Cells.Item($row, 6).Value2 = $newValue
}

# ---------------------------------------------------------------------------
# 2) Add a new "metadata" worksheet, placed right after "data".
# ---------------------------------------------------------------------------
$metaWs = $wb.Worksheets.Add($null, $dataWs)
$metaWs.Name = "metadata"

# Header row (B1:G1), bold/centered like the "data" sheet's header style.
$metaWs.Range("B1").Value2 = "data_name"
$metaWs.Range("C1").Value2 = "data_id"
$metaWs.Range("D1").Value2 = "data_version"
$metaWs.Range("E1").Value2 = "data_version_created"
$metaWs.Range("F1").Value2 = "panel_query_time"
$metaWs.Range("G1").Value2 = "panel_get_request"

# Copy the header formatting (bold font, thin border, centered) from the
# "data" sheet's own header cell so the new tab matches the workbook style.
$dataWs.Range("B1").Copy() | Out-Null
$metaWs.Range("B1:G1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$metaWs.Range("A2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Data row (row 2).
$metaWs.Range("A2").Value2 = 0
$metaWs.Range("B2").Value2 = "Fetal hydrops"
$metaWs.Range("C2").Value2 = 144
$metaWs.Range("D2").NumberFormat = "@"
$metaWs.Range("D2").Value2 = "1.35"
$metaWs.Range("E2").Value2 = "2021-09-14T15:08:11.827403Z"
$metaWs.Range("F2").Value2 = "2021-10-05 14:20:26.038161"
$metaWs.Range("G2").Value2 = "https://panelapp.genomicsengland.co.uk/api/v1/panels/144/?format=json"

$dataWs.Select()
